$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.84"
$ws.Range("D3").Value = "'22.78"
$ws.Range("D4").Value = "'6.146"
$ws.Range("D5").Value = "'0.06071"
$ws.Range("D6").Value = "'6.731"
$ws.Range("D7").Value = "'3.448"
$ws.Range("D8").Value = "'1.358"
$ws.Range("D9").Value = "'0.7968"
$ws.Range("D10").Value = "'0.1576"
$ws.Range("D11").Value = "'0.08021"
$ws.Range("D12").Value = "'0.03344"
$ws.Range("D13").Value = "'0.03094"
$ws.Range("D14").Value = "'0.09303"
$ws.Range("D15").Value = "'3.900"
$ws.Range("D16").Value = "'0.001695"
$ws.Range("D17").Value = "'0.04857"
$ws.Range("D18").Value = "'0.0006144"
$ws.Range("D20").Value = "'0.001101"
$ws.Range("D21").Value = "'0.003383"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.689"
$ws.Range("D24").Value = "'2.265"
$ws.Range("D26").Value = "'0.1235"
$ws.Range("D27").Value = "'0.0003017"
$ws.Range("D40").Value = "'0.04567"
$ws.Range("D41").Value = "'0.007134"
$ws.Range("D42").Value = "'0.003902"
$ws.Range("D43").Value = "'0.1111"
$ws.Range("D45").Value = "'0.002972"
$ws.Range("D46").Value = "'0.00005923"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D48").Value = "'0.7504"
$ws.Range("D49").Value = "'0.07066"
$ws.Range("D50").Value = "'0.00001501"
$ws.Range("D51").Value = "'0.01010"
